# Update the "Förändrad" (changed) date column (C) for rows 2-15
# from 2023-09-06 (45175) to 2023-09-14 (45183).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = 45183
    }
}

